$d = $word.ActiveDocument

# 1) Trim the student ID off the byline: "By Adam Gleichman, A48071742" -> "By Adam Gleichman"
$byline = $d.Content
$found = $byline.Find.Execute(
    "By Adam Gleichman, A48071742", $true, $false, $false, $false, $false,
    $true, 1, $false, "By Adam Gleichman", 2)

if ($found) {
    # After the replace, $byline has collapsed to exactly the new "By Adam
    # Gleichman" text, so $byline.End is the position right after the "n",
    # i.e. the end of the paragraph's text but before its paragraph mark.

    # 2) Move the "_GoBack" bookmark there (it previously sat right after the
    #    final citation field near the end of the document - Word always
    #    re-stamps "_GoBack" at the site of the most recent edit, which is now
    #    this paragraph). Adding a bookmark under a name that already exists
    #    moves it, so the old bookmarkStart/bookmarkEnd pair is removed as a
    #    side effect of this call.
    #
    #    A zero-width Range sitting exactly at "end of paragraph text, before
    #    the paragraph mark" can't be used directly to seed Bookmarks.Add, so
    #    insert a one-character placeholder there first, anchor the bookmark
    #    against that now-interior position, then delete the placeholder -
    #    the (now collapsed) bookmark stays right where the text used to end.
    $tail = $d.Range($byline.End, $byline.End)
    $tail.InsertAfter("#")

    $anchor = $d.Range($byline.End, $byline.End)
    $d.Bookmarks.Add("_GoBack", $anchor) | Out-Null

    $placeholder = $d.Range($byline.End, $byline.End + 1)
    $placeholder.Delete()
}
